# Billboard workbook update ("Added album covers, merged Billboard data"):
#  - Album title "Age of Adz" -> "The Age of Adz"
#  - Header "weeks" -> "weeks_on"
#  - Columns A:C widened to fit the newly merged, longer titles
#  - Active cell moved to H3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the album title in row 4, column A.
$ws.Range("A4").Value = "The Age of Adz"

# Rename the "weeks" header (D1) to "weeks_on".
$ws.Range("D1").Value = "weeks_on"

# Widen columns A:C so the longer merged titles fit.
# (ColumnWidth is rounded to whole on-screen pixels before being stored back
# as a character-width, so nudge the input to land exactly on target.)
function Set-ExactColumnWidth($col, $targetWidth) {
    $ws.Columns.Item($col).ColumnWidth = ([Math]::Round($targetWidth * 6) - 5) / 6
}

Set-ExactColumnWidth 1 16
Set-ExactColumnWidth 2 16.7109375
Set-ExactColumnWidth 3 14.28515625

# Move the active selection to H3.
[void]$ws.Range("H3").Select()
